$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "ID" column header
$ws.Range("C1").Value = "ID"

# Add ID values for existing rows
$ws.Range("C2").Value = 320620321
$ws.Range("C3").Value = 123456789

# Set column C width to fit contents (best fit)
$ws.Columns.Item(3).ColumnWidth = 9.140625

# Update selection to match the new active cell
$ws.Range("C3").Select()
